$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 34, shifting rows 34:67 down to 35:68
# (carries the date number-format on column D down automatically).
$ws.Rows.Item(34).Insert()

# Populate the newly inserted row 34 with the new weekly record.
$ws.Range("A34").Value = 11
$ws.Range("B34").Value = "Vega Monumental Concepción"
$ws.Range("C34").Value = "Bíobío"
$ws.Range("D34").Value = 44589
$ws.Range("E34").Value = 8
$ws.Range("F34").Value = 100112001
$ws.Range("G34").Value = "Berenjena"
$ws.Range("H34").Value = "Sin especificar"
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 180
$ws.Range("K34").Value = 9000
$ws.Range("L34").Value = 10000
$ws.Range("M34").Value = 9556
$ws.Range("N34").Value = '$/caja 60 unidades'
$ws.Range("O34").Value = "Región de Arica y Parinacota"
$ws.Range("P34").Value = 159
$ws.Range("Q34").Value = 60
$ws.Range("R34").Value = "Hortaliza"
